$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @'
63|D.D KENZO/A-008|SALWAR-SUIT|660|D.D KENZO/A-008 - 20/30 - 660|https://i.postimg.cc/k5pYgCs3/Whats-App-Image-2025-05-28-at-16-59-41.jpg
64|D.D KENZO/A-0017|SALWAR-SUIT|660|D.D KENZO/A-0017 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/1ztYSXpr/Whats-App-Image-2025-05-28-at-16-59-42.jpg
65|D.D KENZO/A-0017|SALWAR-SUIT|660|D.D KENZO/A-0017 - 20/36 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/J0wTjQSp/Whats-App-Image-2025-05-28-at-16-59-43.jpg
66|D.D KENZO/A-0016|SALWAR-SUIT|660|D.D KENZO/A-0016 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/VL8G6Dq6/Whats-App-Image-2025-05-28-at-16-59-44.jpg
67|D.D KENZO/A-0016|SALWAR-SUIT|660|D.D KENZO/A-0016 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/dQdNwzTp/Whats-App-Image-2025-05-28-at-16-59-44-1.jpg
68|D.D KENZO/A-0017|SALWAR-SUIT|660|D.D KENZO/A-0017 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/pLVsQr1M/Whats-App-Image-2025-05-28-at-16-59-45.jpg
69|D.D KENZO/A-0018|SALWAR-SUIT|660|D.D KENZO/A-0018 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/tCNrVWkz/Whats-App-Image-2025-05-28-at-16-59-45-1.jpg
70|D.D KENZO/A-001|SALWAR-SUIT|660|D.D KENZO/A-001 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/QMHSPP36/Whats-App-Image-2025-05-28-at-16-59-47.jpg
71|D.D KENZO/A-0016|SALWAR-SUIT|660|D.D KENZO/A-0016 - 20/30 - 3PCS BOX - COTTON - 660|https://i.postimg.cc/3xLnthP4/Whats-App-Image-2025-05-28-at-16-59-47-1.jpg
72|D.D KENZO/A-5366|LACHA|1150|D.D KENZO/A-5366 - 24/34 - 3PCS BOX - 1150|https://i.postimg.cc/SK3VprQm/Whats-App-Image-2025-05-28-at-16-59-49-1.jpg
73|D.D KENZO/A-5353|LACHA|1250|D.D KENZO/A-5353 - 24/34 - 3PCS BOX - 1250|https://i.postimg.cc/d3m9rc20/Whats-App-Image-2025-05-28-at-16-59-50.jpg
74|D.D KENZO/A-5376|LACHA|1190|D.D KENZO/A-5376 - 24/34 - 3PCS BOX - 1190|https://i.postimg.cc/0jqfzyJ3/Whats-App-Image-2025-05-28-at-16-59-50-1.jpg
75|D.D KENZO/A-6001|FARSHI|1190|D.D KENZO/A-6001 - 24/34 - 3PCS BOX - 1190|https://i.postimg.cc/zX8753NH/Whats-App-Image-2025-05-28-at-16-59-51.jpg
76|D.D. KENZO/A-5375|LACHA|1090|D.D. KENZO/A-5375 - 24/34 - 3PCS BOX - 1090|https://i.postimg.cc/FFjp23MQ/Whats-App-Image-2025-05-28-at-16-59-51-1.jpg
77|D.D KENZO/A-5360|LACHA|1290|D.D KENZO/A-5360 - 24/34 - 3PCS BOX - 1290|https://i.postimg.cc/X7L8v3NH/Whats-App-Image-2025-05-28-at-16-59-53.jpg
78|D.D KENZO/A-5367|LACHA|1280|D.D KENZO/A-5367 - 24/34 - 3PCS BOX - 1280|https://i.postimg.cc/T3HJr9Dn/Whats-App-Image-2025-05-28-at-16-59-54.jpg
79|D.D KENZO/A-6102|LACHA|710|D.D KENZO/A-6102 - 18/22 - 3PCS BOX - 710|https://i.postimg.cc/g20KX5nY/Whats-App-Image-2025-05-28-at-16-59-54-1.jpg
80|D.D KENZO/A-6076|LACHA|750|D.D KENZO/A-6076 - 18/22 - 3PCS BOX - 750|https://i.postimg.cc/cJgMTnR2/Whats-App-Image-2025-05-28-at-16-59-55.jpg
81|CUTIES/A-4353|SHARARA|1120|CUTIES/A-4353 - 24/34 - 3PCS BOX - 1120|https://i.postimg.cc/ZR6PzzDL/Whats-App-Image-2025-05-28-at-16-59-56.jpg
82|D.D /A-6106|LACHA|730|D.D /A-6106 - 18/22 - 3PCS BOX - 730|https://i.postimg.cc/fbG7PR5L/Whats-App-Image-2025-05-28-at-16-59-56-1.jpg
83|CUTIES/A-4336|DHOTI|1120|CUTIES/A-4336 - 24/34 - 3PCS BOX - 1120|https://i.postimg.cc/tJ5h7tbZ/Whats-App-Image-2025-05-28-at-16-59-57.jpg
84|CUTIES/A-4355|SHARARA|1120|CUTIES/A-4355 - 24/34 - 3PCS BOX - 1120|https://i.postimg.cc/7hD3xtFr/Whats-App-Image-2025-05-28-at-16-59-57-1.jpg
85|D.D. KENZO|SHARARA|1290|D.D. KENZO - 24/34 - 3PCS BOX - 1290|https://i.postimg.cc/28kQMJd4/Whats-App-Image-2025-05-28-at-16-59-58.jpg
86|D.D KENZO|SHARARA|1260|D.D KENZO - 24/34 - 3PCS BOX - 1260|https://i.postimg.cc/ryCGhc4V/Whats-App-Image-2025-05-28-at-16-59-58-1.jpg
87|D.D KENZO/A-483|GHARARA|1120|D.D KENZO/A-483 - 24/34 - 3PCS BOX - 1120|https://i.postimg.cc/3xWgJvcN/Whats-App-Image-2025-05-28-at-16-59-59.jpg
88|D.D KENZO/A-468|GHARARA|1140|D.D KENZO/A-468 - 24/34 - 3PCS BOX - 1140|https://i.postimg.cc/zDdD1zyD/Whats-App-Image-2025-05-29-at-13-24-50.jpg
89|D.D/A-9132|PLAZO|690|D.D/A-9132 - 18/22 - 3PCS BOX - 690|https://i.postimg.cc/3R8hPV7D/Whats-App-Image-2025-05-29-at-13-24-44.jpg
90|D.D/A-9133|GHARARA|720|D.D/A-9133 - 18/22 - 3PCS BOX - 720|https://i.postimg.cc/hj2ngsPZ/Whats-App-Image-2025-05-29-at-13-24-45.jpg
91|D.D/A-9124|GHARARA|690|D.D/A-9124 - 18/22 - 3PCS BOX - 690|https://i.postimg.cc/5yGMQtzw/Whats-App-Image-2025-05-29-at-13-24-45-1.jpg
92|D.D/A-9110|GHARARA|655|D.D/A-9110 - 18/22 - 3PCS BOX - 655|https://i.postimg.cc/kX6C3xrw/Whats-App-Image-2025-05-29-at-13-24-46.jpg
93|PINK CASTLE/A-1180|CAPRI SUIT|595|PINK CASTLE/A-1180 - 16/32 - 3PCS BOX - COTTON & RFD - 595|https://i.postimg.cc/Bbt0vFjS/Whats-App-Image-2025-05-29-at-13-24-46-1.jpg
94|TINA MINA/A-1250|SKIRT-TOP|715|TINA MINA/A-1250 - 22/32 - 3PCS BOX - 715|https://i.postimg.cc/WbBv2hTV/Whats-App-Image-2025-05-29-at-13-24-47.jpg
95|TINA MINA/A-1259|NICKER-SUIT|650|TINA MINA/A-1259 - 22/32 - 3PCS BOX - 650|https://i.postimg.cc/RCbmDgYN/Whats-App-Image-2025-05-29-at-13-24-47-1.jpg
96|MINI-Q A/7547|GIRLISH SUIT|350|MINI-Q A/7547 - 22/38 - 3PCS BOX - 350|https://i.postimg.cc/cJmyMhY1/Whats-App-Image-2025-05-29-at-13-24-48.jpg
97|MINI-Q A/-7541|GIRLISH SUIT|330|MINI-Q A/-7541 - 22/38 - 3PCS BOX - COTTON LYCRA - 330|https://i.postimg.cc/JzjWtMzk/Whats-App-Image-2025-05-29-at-13-24-48-1.jpg
98|SALVIA/A-1543|CORD-SET|335|SALVIA/A-1543 - 22/38 - 3PCS BOX - COTTON LYCRA - 335|https://i.postimg.cc/9M86sfJv/Whats-App-Image-2025-05-29-at-13-24-49.jpg
99|SALVIA/A-1531|CORD-SET|335|SALVIA/A-1531 - 22/32 - 3PCS BOX - COTTON LYCRA - 335|https://i.postimg.cc/Y9CcXFf6/Whats-App-Image-2025-05-29-at-13-24-49-1.jpg
100|SMILE /A-2135|CORD-SET|445|SMILE /A-2135 - 16/32 - 3PCS BOX - COTTON LYCRA - 445|https://i.postimg.cc/Kj16byW4/Whats-App-Image-2025-05-29-at-13-24-50-1.jpg
101|MINNIE /A-1130|NICKER SUIT|415|MINNIE /A-1130 - 16/26 - 3PCS BOX - 415|https://i.postimg.cc/02JTHPmB/Whats-App-Image-2025-05-29-at-13-24-51.jpg
102|TINI MINI/A-7138|CORD-SET|285|TINI MINI/A-7138 - 22/32 - 3PCS BOX - TENCIL COTTON - 285|https://i.postimg.cc/HWKG2VJM/Whats-App-Image-2025-05-29-at-13-24-52.jpg
103|TINI MINI/A-7233|CAPRI SUIT|205|TINI MINI/A-7233 - 22/38 - 3PCS BOX - COTTON - 205|https://i.postimg.cc/gj1fkT9s/Whats-App-Image-2025-05-29-at-13-24-52-1.jpg
104|TINI MINI/A-7136|CAPRI SUIT|205|TINI MINI/A-7136 - 22/32 - 3PCS BOX - TENCIL COTTON - 205|https://i.postimg.cc/Vkq3mgK9/Whats-App-Image-2025-05-29-at-13-24-53.jpg
105|TINI MINI /A-7232|CAPRI SUIT|165|TINI MINI /A-7232 - 22/32 - 3PCS BOX - COTTON - 165|https://i.postimg.cc/CxTW3XsL/Whats-App-Image-2025-05-29-at-13-24-53-1.jpg
106|TINI MINI/A-7210|NICKER SUIT|240|TINI MINI/A-7210 - 22/32 - 3PCS BOX - COTTON - 240|https://i.postimg.cc/d0Jg1Cfb/Whats-App-Image-2025-05-29-at-13-24-54.jpg
107|BALENO/A-0000|T-SHIRT|97|BALENO/A-0000 - 20/36 - 18PCS BOX - 97|https://i.postimg.cc/1tPgvwB5/Whats-App-Image-2025-05-29-at-13-28-56.jpg
108|GLARE/A-106|T-SHIRT|108|GLARE/A-106 - 18/34 - 18PCS BOX - COTTON - 108|https://i.postimg.cc/6qby6SXy/Whats-App-Image-2025-05-29-at-13-28-57.jpg
109|SPARK/A-7546|T-SHIRT|377|SPARK/A-7546 - 26/36 - 3PCS BOX - COTTON - 377|https://i.postimg.cc/rmJD3sZG/Whats-App-Image-2025-05-29-at-13-28-57-1.jpg
110|SPARK/A-7536|T-SHIRT|405|SPARK/A-7536 - 26/36 - 3PCS BOX - COTTON - 405|https://i.postimg.cc/7hZ5k4jF/Whats-App-Image-2025-05-29-at-13-28-58.jpg
111|BALENO/A-7408|T-SHIRT|125|BALENO/A-0000 - 6PCS BOX - 97|https://i.postimg.cc/25nbd24K/Whats-App-Image-2025-05-29-at-13-28-58-1.jpg
112|BALENO/A-7410|T-SHIRT|125|BALENO/A-7410 - 20/36 - 6PCS BOX - COTTON (D-SHOULDER) - 125|https://i.postimg.cc/bNdZVzxb/Whats-App-Image-2025-05-29-at-13-29-04-1.jpg
113|BALENO/A-6407|T-SHIRT|195|BALENO/A-6407 - 20/36 - 3PCS BOX - COTTON - 195|https://i.postimg.cc/YqVhJJRP/Whats-App-Image-2025-05-29-at-13-29-05.jpg
114|BALENO/A-6405|T-SHIRT|165|BALENO/A-6405 - 20/36 - 3PCS BOX - TENCIL - 165|https://i.postimg.cc/wBY3BLSc/Whats-App-Image-2025-05-29-at-13-29-05-1.jpg
115|BALENO/A-6409|T-SHIRT|120|BALENO/A-6409 - 20/36 - 3PCS BOX - TENCIL - 120|https://i.postimg.cc/fLwLcZWV/Whats-App-Image-2025-05-29-at-13-29-06.jpg
116|BALENO/A-7853|T-SHIRT|195|BALENO/A-7853 - 20/36 - 3PCS BOX - TENCIL - 195|https://i.postimg.cc/fb7ySks4/Whats-App-Image-2025-05-29-at-13-29-06-1.jpg
117|CRISPY/A-2330|KURTA-PAJAMA|515|CRISPY/A-2330 - 2/9 - 1PCS BOX - 515|https://i.postimg.cc/VLqgKWm2/Whats-App-Image-2025-05-29-at-15-38-23.jpg
118|CRISPY/A-2330|KURTA-PAJAMA|515|CRISPY/A-2330 - 2/9 - 1PCS BOX - 515|https://i.postimg.cc/gkF4yTFv/Whats-App-Image-2025-05-29-at-15-38-24.jpg
119|NEW TOPS/A-6080|KURTA-PAJAMA|485|NEW TOPS/A-6080 - 20/36 - 1PCS BOX - 485|https://i.postimg.cc/QdjbwKMC/Whats-App-Image-2025-05-29-at-15-38-24-2.jpg
120|FROGGER/A-1487|T-SHIRT|122|FROGGER/A-1487 - 24/36 - 6PCS BOX - TENCIL (D-SHOULDER) - 122|https://i.postimg.cc/c6y95d9S/Whats-App-Image-2025-05-29-at-16-00-16.jpg
121|A /A-1417|T-SHIRT|133|A /A-1417 - 24/36 - 6PCS BOX - TENCIL LYCRA - 133|https://i.postimg.cc/4Nps4yTp/Whats-App-Image-2025-05-29-at-16-00-17.jpg
122|A /A-32391|T-SHIRT|158|A /A-32391 - 26/36 - 3PCS BOX - COTTON (D-SHOULDER) - 158|https://i.postimg.cc/h4zKwMyy/Whats-App-Image-2025-05-29-at-16-00-17-1.jpg
123|A+ /A-32706|T-SHIRT|194|A+ /A-32706 - 26/36 - 3PCS BOX - TENCIL (D-SHOULDER) - 194|https://i.postimg.cc/T3GGKxRF/Whats-App-Image-2025-05-29-at-16-00-18.jpg
124|A+ /A-32569|T-SHIRT|228|A+ /A-32569 - 26/36 - 3PCS BOX - COTTON LYCRA - 228|https://i.postimg.cc/3wxTZZ0Z/Whats-App-Image-2025-05-29-at-16-00-19.jpg
125|A+ /A-2102|T-SHIRT|210|A+ /A-2102 - 24/36 - 6PCS BOX - 210|https://i.postimg.cc/XY8WFPqJ/Whats-App-Image-2025-05-29-at-16-00-19-1.jpg
126|FRESH SPORTS /A-2124|T-SHIRT|210|FRESH SPORTS /A-2124 - 24/36 - 6PCS BOX - 210|https://i.postimg.cc/hjHRcWyC/Whats-App-Image-2025-05-29-at-16-00-20.jpg
127|A /A-B01|T-SHIRT|288|A /A-B01 - 26/36 - 3PCS BOX - IMPORTEED - 288|https://i.postimg.cc/zf212RnY/Whats-App-Image-2025-05-29-at-16-00-20-1.jpg
128|A /A-V01|T-SHIRT|288|A /A-V01 - 26/36 - 3PCS BOX - IMPORTEED - 288|https://i.postimg.cc/GtGwp76b/Whats-App-Image-2025-05-29-at-16-00-21.jpg
129|ALAG /A-969|T-SHIRT|319|ALAG /A-969 - 26/36 - 3PCS BOX - IMPORTEED - 319|https://i.postimg.cc/KzKh6vNk/Whats-App-Image-2025-05-29-at-16-00-21-1.jpg
130|FRESH SPORTS /A-1268|T-SHIRT|425|FRESH SPORTS /A-1268 - 26/36 - 3PCS BOX - COTTON  - 425|https://i.postimg.cc/Prnn7H8J/Whats-App-Image-2025-05-29-at-16-00-22.jpg
131|ALAG/A-904|T-SHIRT|268|ALAG/A-904 - 26/36 - 3PCS BOX - TEXTURE - 268|https://i.postimg.cc/v8VdXTTC/Whats-App-Image-2025-05-29-at-16-00-22-1.jpg
132|FRESH SPORTS/A-1295|T-SHIRT|480|FRESH SPORTS/A-1295 - 26/36 - 3PCS BOX - C WASH - 480|https://i.postimg.cc/vZ3RD62x/Whats-App-Image-2025-05-29-at-16-00-23.jpg
133|FRESH SPORTS/A-1291|T-SHIRT|380|FRESH SPORTS/A-1291 - 26/36 - 3PCS BOX - 380|https://i.postimg.cc/BQ5f9zxg/Whats-App-Image-2025-05-29-at-16-00-23-1.jpg
134|FRESH SPORTS/A-1293|T-SHIRT|395|FRESH SPORTS/A-1293 - 32/36 - 3PCS BOX - MINI WAFFLE - 395|https://i.postimg.cc/k540SXL1/Whats-App-Image-2025-05-29-at-16-00-24.jpg
'@

$lines = $rowData -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]
    $aVal = $parts[1]
    $bVal = $parts[2]
    $cVal = [double]$parts[3]
    $dVal = $parts[4]
    $eVal = $parts[5]

    $ws.Cells.Item($r, 1).Value = $aVal
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 3).HorizontalAlignment = 1
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
}

$ws.Range("D18").Select()

Write-Output "done"
